{"js": "// Rename the document title from \"IDG PA28X\" to \"IDG PA28\" and move the\n// \"_GoBack\" bookmark so it sits right after the (now shorter) title text,\n// matching Word's own behavior of relocating the last-edit-position\n// bookmark to the most recently edited spot.\n\nconst body = context.document.body;\n\n// 1) Find the exact title run text and trim the trailing \"X\".\nconst titleResults = body.search(\"IDG PA28X\", { matchCase: true, matchWholeWord: false });\ntitleResults.load(\"items\");\nawait context.sync();\n\nif (titleResults.items.length === 0) {\n  throw new Error(\"Could not find 'IDG PA28X' text to rename.\");\n}\n\nconst titleRange = titleResults.items[0];\ntitleRange.insertText(\"IDG PA28\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) Remove the existing \"_GoBack\" bookmark (wherever Word last left it).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 3) Re-create \"_GoBack\" immediately after the title text we just edited.\nconst afterTitle = titleRange.getRange(Word.RangeLocation.end);\nafterTitle.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the current title text \"IDG PA28X\".\n$find = $d.Content.Find\n$find.Text = \"IDG PA28X\"\n$found = $find.Execute()\n\nif ($found) {\n  $titleRange = $find.Parent.Duplicate\n\n  # Re-point the \"_GoBack\" bookmark to the end of the title range *before*\n  # trimming the trailing \"X\" - this also removes/relocates any pre-existing\n  # \"_GoBack\" bookmark elsewhere in the document (a document can only have\n  # one bookmark with a given name).\n  $endRange = $titleRange.Duplicate\n  $endRange.Collapse(0) | Out-Null  # wdCollapseEnd\n  $d.Bookmarks.Add(\"_GoBack\", $endRange) | Out-Null\n\n  # Remove the trailing \"X\" from \"IDG PA28X\" -> \"IDG PA28\".\n  $delRange = $titleRange.Duplicate\n  $delRange.Collapse(0) | Out-Null   # wdCollapseEnd\n  $delRange.MoveStart(1, -1) | Out-Null  # wdCharacter: shrink start back by 1 char\n  $delRange.Delete() | Out-Null\n}\n"}
